$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Retitle cm019 and cm020 topics
$ws.Range("D20").Value = "Building Shiny applications: user interface"
$ws.Range("D21").Value = "Building Shiny applications: server"
